$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.36823296546936
$ws.Range("B1").Value = 1.90549623966217
$ws.Range("C1").Value = 3.385767698287964
$ws.Range("D1").Value = 3.752682209014893
$ws.Range("E1").Value = 0.9806156158447266
